# Updated symbol list on Sat Dec 17 02:55:19 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) values for a
# number of coins on Sheet1, matching the latest scrape of the symbol list.
# Price values are stored as text (inline strings) in the workbook, so each
# one is written with a Text number format to stop Excel from re-interpreting
# the numeric-looking string as a number, then the style is restored back to
# Normal so no visible formatting change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Column D - Price
Set-TextValue "D2"  "227.52"
Set-TextValue "D3"  "22.54"
Set-TextValue "D4"  "5.330"
Set-TextValue "D5"  "0.05527"
Set-TextValue "D6"  "3.388"
Set-TextValue "D8"  "0.7830"
Set-TextValue "D9"  "1.045"
Set-TextValue "D10" "0.1384"
Set-TextValue "D11" "0.07421"
Set-TextValue "D12" "0.03138"
Set-TextValue "D13" "0.02939"
Set-TextValue "D14" "0.09253"
Set-TextValue "D15" "0.001670"
Set-TextValue "D16" "3.254"
Set-TextValue "D17" "0.04768"
Set-TextValue "D18" "0.0005861"
Set-TextValue "D19" "0.006217"
Set-TextValue "D20" "0.005229"
Set-TextValue "D21" "0.001064"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.862"
Set-TextValue "D24" "2.197"
Set-TextValue "D26" "0.1283"
Set-TextValue "D27" "0.0005021"
Set-TextValue "D40" "0.03906"
Set-TextValue "D41" "0.007134"
Set-TextValue "D42" "0.1033"
Set-TextValue "D43" "0.003267"
Set-TextValue "D44" "0.009239"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.7853"
Set-TextValue "D48" "0.08873"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.01010"

# Column E - Volume(1h) (label such as "Best/Worst in 24h" toggled for a few coins)
$ws.Range("E9").Value  = "8FTXTokenFTTWorstin24h"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
